# ---------------------------------------------------------------------------
# "杭州-漫展信息.xlsx" data refresh
#
# The site was re-scraped and most events simply gained a few more
# "want to go" (column F) counters. On top of that, sheet "全部类型"
# (the combined/union sheet) lost its "天空之城" concert row and gained the
# "第三届百合Only" expo row that already existed on sheet "展览" - which
# shifts every row between them down... sorry, UP by one, filling rows
# 34-39 with what used to be rows 35-40, and inserting the new row at 40.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$wsExhibit = $wb.Worksheets.Item(1)   # 展览   (exhibitions)
$wsShow    = $wb.Worksheets.Item(2)   # 演出   (performances)
$wsLocal   = $wb.Worksheets.Item(3)   # 本地生活 (untouched)
$wsAll     = $wb.Worksheets.Item(4)   # 全部类型 (union of all of the above)

# ---------------------------------------------------------------------------
# 1) Plain "想去人数" (want-to-go count, column F) refreshes on 展览.
#    Nothing else moves on this sheet.
# ---------------------------------------------------------------------------
$exhibitUpdates = @{
    2  = 1123
    3  = 255
    4  = 260
    6  = 676
    7  = 535
    8  = 4906
    10 = 453
    11 = 115
    12 = 992
    13 = 333
    14 = 1296
    18 = 1839
    20 = 49
    21 = 172
    22 = 69
    24 = 946
    26 = 35
    27 = 3086
    28 = 1040
    29 = 2526
    30 = 256
    31 = 1381
    32 = 3693
    33 = 97
    34 = 899
    36 = 1160
    37 = 6
    38 = 940
    39 = 1203
    40 = 28
    41 = 899
    42 = 579
    43 = 327
    44 = 372
    46 = 3504
}
foreach ($row in $exhibitUpdates.Keys) {
    $wsExhibit.Cells.Item([int]$row, 6).Value = $exhibitUpdates[$row]
}

# ---------------------------------------------------------------------------
# 2) Same kind of refresh on 演出.
# ---------------------------------------------------------------------------
$showUpdates = @{
    5  = 3
    11 = 889
    19 = 2
    23 = 33
}
foreach ($row in $showUpdates.Keys) {
    $wsShow.Cells.Item([int]$row, 6).Value = $showUpdates[$row]
}

# ---------------------------------------------------------------------------
# 3) 全部类型: the rows that only need their column F refreshed in place
#    (everything above row 34 and the tail rows 41+ after the shuffle).
# ---------------------------------------------------------------------------
$allTypeUpdates = @{
    2  = 1123
    3  = 255
    4  = 260
    7  = 676
    8  = 535
    9  = 4906
    13 = 333
    14 = 1296
    17 = 1839
    19 = 49
    21 = 172
    22 = 889
    25 = 69
    26 = 946
    28 = 3086
    30 = 1040
    31 = 2526
    32 = 1381
    33 = 3693
    41 = 900
    42 = 579
    43 = 372
    44 = 33
    48 = 3504
}
foreach ($row in $allTypeUpdates.Keys) {
    $wsAll.Cells.Item([int]$row, 6).Value = $allTypeUpdates[$row]
}

# ---------------------------------------------------------------------------
# 4) 全部类型 rows 34-40: the "天空之城" row disappears and "第三届百合Only"
#    is pulled in from 展览, so every row in between shifts up by one.
#    Re-source each row's B-I content fresh off 展览 / 演出 (column A, the
#    running index, stays put - it never changed in the source diff).
# ---------------------------------------------------------------------------
function Copy-RowToAllTypes($srcWs, $srcRow, $dstRow) {
    # Column B holds a plain "YYYY-MM-DD" label. Assigning that text as-is
    # gets auto-promoted to a real date value, so force it to stay text the
    # same way a user would in Excel (leading apostrophe), then drop the
    # resulting quote-prefix formatting so the cell style is untouched.
    $bVal = $srcWs.Cells.Item($srcRow, 2).Value2
    $dst = $wsAll.Cells.Item($dstRow, 2)
    $dst.Value = "'" + $bVal
    $dst.Style = "Normal"

    foreach ($col in 3..9) {
        $dst = $wsAll.Cells.Item($dstRow, $col)
        $dst.Value = $srcWs.Cells.Item($srcRow, $col).Value2
        $dst.Style = "Normal"
    }
}

Copy-RowToAllTypes $wsExhibit 33 34   # 杭州·亚米二次茶话会展
Copy-RowToAllTypes $wsExhibit 34 35   # 杭州·次元幻想--二次元全女夜场
Copy-RowToAllTypes $wsExhibit 36 36   # 杭州·第五届华盟次元嘉年华&周年庆狂欢
Copy-RowToAllTypes $wsExhibit 38 37   # 杭州·夏之誓国乙only-日夜场
Copy-RowToAllTypes $wsShow    20 38   # 杭州·黄西全新脱口秀专场《水土不服》
Copy-RowToAllTypes $wsExhibit 39 39   # 杭州·火影忍者only
Copy-RowToAllTypes $wsExhibit 40 40   # 杭州·第三届百合Only·同好交流 (new row)

# These titles picked up extra "want to go" votes between the two scrapes,
# same as everywhere else - apply on top of the freshly-copied rows.
$wsAll.Cells.Item(34, 6).Value = 97
$wsAll.Cells.Item(35, 6).Value = 899
$wsAll.Cells.Item(36, 6).Value = 1160
$wsAll.Cells.Item(37, 6).Value = 940
$wsAll.Cells.Item(39, 6).Value = 1203
$wsAll.Cells.Item(40, 6).Value = 28
